$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a numeric-looking value while keeping the cell's underlying
# type as TEXT (the sheet stores all prices as inline strings, not numbers).
# We briefly flip the cell to Text number format so Excel doesn't coerce the
# numeric-looking string into a real number, then restore the original
# "Normal" style so no stray formatting is left behind.
function Set-TextValue($addr, $value) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

# Plain text updates (no numeric coercion risk).
function Set-TextCell($addr, $value) {
    $ws.Range($addr).Value = $value
}

# --- Price (column D) updates for unchanged rows ---
Set-TextValue "D2"  "245.30"
Set-TextValue "D3"  "22.06"
Set-TextValue "D4"  "5.331"
Set-TextValue "D5"  "0.05974"
Set-TextValue "D6"  "3.402"
Set-TextValue "D7"  "6.385"
Set-TextValue "D8"  "0.8113"
Set-TextValue "D9"  "0.9682"
Set-TextValue "D10" "0.1426"

# --- Rows 11 and 12 swap places (coin name / link / volume label), with new prices ---
Set-TextCell  "B11" "MandalaExchangeToken"
Set-TextCell  "C11" "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D11" "0.07415"
Set-TextCell  "E11" "10MandalaExchangeTokenMDX"

Set-TextCell  "B12" "LiechtensteinCryptoassetsExchange"
Set-TextCell  "C12" "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue "D12" "0.03500"
Set-TextCell  "E12" "11LiechtensteinCryptoassetsExchangeLCX"

# --- Remaining price updates ---
Set-TextValue "D13" "0.03071"
Set-TextValue "D15" "3.997"
Set-TextValue "D16" "0.001600"
Set-TextValue "D17" "0.04795"

Set-TextCell  "E18" "17OneONEWorstin24h"

Set-TextValue "D19" "0.006236"
Set-TextValue "D20" "0.004137"
Set-TextValue "D21" "0.0009891"
Set-TextValue "D22" "0.00009708"
Set-TextValue "D24" "2.165"
Set-TextValue "D26" "0.1333"
Set-TextValue "D40" "0.03911"
Set-TextValue "D41" "0.006496"
Set-TextValue "D42" "0.1073"
Set-TextValue "D43" "0.002702"
Set-TextValue "D44" "0.005381"
Set-TextValue "D45" "0.00005321"
Set-TextValue "D48" "0.03998"

Set-TextCell  "E48" "47BOLOBOLO"
